# Minor updates to dataset
# - Fix dog name typo "Cozzy" -> "Cozy" (row 11)
# - Fix Neutered/Sprayed value for Diesel (row 15) from "Unknown" -> "N/A"
# - Update the active cell selection to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "Cozy"
$ws.Range("F15").Value = "N/A"

[void]$ws.Range("D13").Select()
